$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label for column B from "unnamed: 1_level_1" to "total"
$ws.Range("B2").Value = "total"

# Update the data table values (rows 4-38, columns B-H)
$ws.Range("B4").Value = 1.44
$ws.Range("C4").Value = 4.87
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 2.18
$ws.Range("F4").Value = 1.99
$ws.Range("G4").Value = 2.19
$ws.Range("H4").Value = 9.09
$ws.Range("B5").Value = 1.43
$ws.Range("C5").Value = 4.99
$ws.Range("D5").Value = 3.02
$ws.Range("E5").Value = 2.22
$ws.Range("F5").Value = 2.04
$ws.Range("G5").Value = 2.21
$ws.Range("H5").Value = 9.74
$ws.Range("B6").Value = 5.49
$ws.Range("C6").Value = 10.61
$ws.Range("D6").Value = 9.41
$ws.Range("E6").Value = 7.76
$ws.Range("F6").Value = 8.41
$ws.Range("G6").Value = 10.76
$ws.Range("H6").Value = 21.94
$ws.Range("B7").Value = 4.95
$ws.Range("C7").Value = 12.44
$ws.Range("D7").Value = 7.53
$ws.Range("E7").Value = 5.92
$ws.Range("F7").Value = 5.22
$ws.Range("G7").Value = 5.49
$ws.Range("H7").Value = 21.33
$ws.Range("B8").Value = 8.82
$ws.Range("C8").Value = 32.51
$ws.Range("D8").Value = 18.32
$ws.Range("E8").Value = 10.35
$ws.Range("F8").Value = 12.5
$ws.Range("G8").Value = 13.7
$ws.Range("H8").Value = 52.65
$ws.Range("B9").Value = 15.85
$ws.Range("C9").Value = 37.09
$ws.Range("D9").Value = 29.21
$ws.Range("E9").Value = 18.97
$ws.Range("F9").Value = 16.2
$ws.Range("G9").Value = 13.11
$ws.Range("H9").Value = 55.05
$ws.Range("B10").Value = 11.66
$ws.Range("C10").Value = 24.92
$ws.Range("D10").Value = 23.03
$ws.Range("E10").Value = 17
$ws.Range("F10").Value = 14.36
$ws.Range("G10").Value = 18.26
$ws.Range("H10").Value = 50.71
$ws.Range("B11").Value = 11.4
$ws.Range("C11").Value = 32.4
$ws.Range("D11").Value = 23.63
$ws.Range("E11").Value = 19.78
$ws.Range("F11").Value = 17.47
$ws.Range("G11").Value = 16.54
$ws.Range("H11").Value = 95.34999999999999
$ws.Range("B12").Value = 7.66
$ws.Range("C12").Value = 17.48
$ws.Range("D12").Value = 10.06
$ws.Range("E12").Value = 8.550000000000001
$ws.Range("F12").Value = 8.16
$ws.Range("G12").Value = 8.59
$ws.Range("H12").Value = 32.51
$ws.Range("B13").Value = 16.45
$ws.Range("C13").Value = 38.22
$ws.Range("D13").Value = 24.51
$ws.Range("E13").Value = 17.88
$ws.Range("F13").Value = 19.69
$ws.Range("G13").Value = 21.37
$ws.Range("H13").Value = 56.97
$ws.Range("B14").Value = 13.73
$ws.Range("C14").Value = 24.56
$ws.Range("D14").Value = 29.52
$ws.Range("E14").Value = 14.85
$ws.Range("F14").Value = 14.05
$ws.Range("G14").Value = 12.8
$ws.Range("H14").ClearContents()
$ws.Range("B15").Value = 2.68
$ws.Range("C15").Value = 6.77
$ws.Range("D15").Value = 4.67
$ws.Range("E15").Value = 3.58
$ws.Range("F15").Value = 3.69
$ws.Range("G15").Value = 4.13
$ws.Range("H15").Value = 14.63
$ws.Range("B16").Value = 14.88
$ws.Range("C16").Value = 33.57
$ws.Range("D16").Value = 23.98
$ws.Range("E16").Value = 15.07
$ws.Range("F16").Value = 12.28
$ws.Range("G16").Value = 19.66
$ws.Range("H16").Value = 28.54
$ws.Range("B17").Value = 14.69
$ws.Range("C17").Value = 33
$ws.Range("D17").Value = 29.72
$ws.Range("E17").Value = 20.7
$ws.Range("F17").Value = 16.68
$ws.Range("G17").Value = 20.36
$ws.Range("H17").Value = 96.81999999999999
$ws.Range("B18").Value = 4.51
$ws.Range("C18").Value = 10.57
$ws.Range("D18").Value = 8.23
$ws.Range("E18").Value = 7.33
$ws.Range("F18").Value = 6.47
$ws.Range("G18").Value = 8.720000000000001
$ws.Range("H18").Value = 45.15
$ws.Range("B19").Value = 7.28
$ws.Range("C19").Value = 27.77
$ws.Range("D19").Value = 15.65
$ws.Range("E19").Value = 10.82
$ws.Range("F19").Value = 12.71
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 52.36
$ws.Range("B20").Value = 9.789999999999999
$ws.Range("C20").Value = 23.88
$ws.Range("D20").Value = 13.7
$ws.Range("E20").Value = 16.53
$ws.Range("F20").Value = 15.19
$ws.Range("G20").Value = 16.82
$ws.Range("H20").Value = 67.26000000000001
$ws.Range("B21").Value = 5.35
$ws.Range("C21").Value = 12.22
$ws.Range("D21").Value = 9.470000000000001
$ws.Range("E21").Value = 7.93
$ws.Range("F21").Value = 8.699999999999999
$ws.Range("G21").Value = 9.039999999999999
$ws.Range("H21").Value = 26.78
$ws.Range("B22").Value = 12.79
$ws.Range("C22").Value = 27.13
$ws.Range("D22").Value = 18.02
$ws.Range("E22").Value = 18.58
$ws.Range("F22").Value = 24.87
$ws.Range("G22").Value = 23.85
$ws.Range("H22").Value = 55.55
$ws.Range("B23").Value = 7.29
$ws.Range("C23").Value = 25.8
$ws.Range("D23").Value = 12.4
$ws.Range("E23").Value = 11.81
$ws.Range("F23").Value = 10.53
$ws.Range("G23").Value = 15.73
$ws.Range("H23").Value = 89.44
$ws.Range("B24").Value = 5.24
$ws.Range("C24").Value = 12.97
$ws.Range("D24").Value = 8.42
$ws.Range("E24").Value = 7.56
$ws.Range("F24").Value = 7.73
$ws.Range("G24").Value = 6.22
$ws.Range("H24").Value = 30.65
$ws.Range("B25").Value = 2.51
$ws.Range("C25").Value = 11.76
$ws.Range("D25").Value = 6.53
$ws.Range("E25").Value = 4.38
$ws.Range("F25").Value = 3.49
$ws.Range("G25").Value = 3.54
$ws.Range("H25").Value = 17.46
$ws.Range("B26").Value = 4.69
$ws.Range("C26").Value = 20.99
$ws.Range("D26").Value = 11.35
$ws.Range("E26").Value = 7.7
$ws.Range("F26").Value = 5.95
$ws.Range("G26").Value = 6.8
$ws.Range("H26").Value = 35.51
$ws.Range("B27").Value = 8.779999999999999
$ws.Range("C27").Value = 28.68
$ws.Range("D27").Value = 19.45
$ws.Range("E27").Value = 14.2
$ws.Range("F27").Value = 14.14
$ws.Range("G27").Value = 13.1
$ws.Range("H27").Value = 68.73999999999999
$ws.Range("B28").Value = 4.47
$ws.Range("C28").Value = 27.05
$ws.Range("D28").Value = 13.89
$ws.Range("E28").Value = 7.97
$ws.Range("F28").Value = 6.74
$ws.Range("G28").Value = 6.21
$ws.Range("H28").Value = 30.56
$ws.Range("B29").Value = 3.72
$ws.Range("C29").Value = 17.28
$ws.Range("D29").Value = 10.33
$ws.Range("E29").Value = 6.98
$ws.Range("F29").Value = 5.08
$ws.Range("G29").Value = 5.06
$ws.Range("H29").Value = 26.05
$ws.Range("B30").Value = 3.31
$ws.Range("C30").Value = 15
$ws.Range("D30").Value = 8.949999999999999
$ws.Range("E30").Value = 5.63
$ws.Range("F30").Value = 4.52
$ws.Range("G30").Value = 4.55
$ws.Range("H30").Value = 25.17
$ws.Range("B31").Value = 4.47
$ws.Range("C31").Value = 22.55
$ws.Range("D31").Value = 13.28
$ws.Range("E31").Value = 8.23
$ws.Range("F31").Value = 6.3
$ws.Range("G31").Value = 6.76
$ws.Range("H31").Value = 36.78
$ws.Range("B32").Value = 11.43
$ws.Range("C32").Value = 51.2
$ws.Range("D32").Value = 27.48
$ws.Range("E32").Value = 14.41
$ws.Range("F32").Value = 14.58
$ws.Range("G32").Value = 13
$ws.Range("H32").Value = 74.13
$ws.Range("B33").Value = 4.77
$ws.Range("C33").Value = 20.77
$ws.Range("D33").Value = 13.42
$ws.Range("E33").Value = 9.199999999999999
$ws.Range("F33").Value = 6
$ws.Range("G33").Value = 5.94
$ws.Range("H33").Value = 35.78
$ws.Range("B34").Value = 3.36
$ws.Range("C34").Value = 16.14
$ws.Range("D34").Value = 7.65
$ws.Range("E34").Value = 5.09
$ws.Range("F34").Value = 4.78
$ws.Range("G34").Value = 4.96
$ws.Range("H34").Value = 24.29
$ws.Range("B35").Value = 15.59
$ws.Range("C35").Value = 54.6
$ws.Range("D35").Value = 24.34
$ws.Range("E35").Value = 15.64
$ws.Range("F35").Value = 17.84
$ws.Range("G35").Value = 16.7
$ws.Range("H35").Value = 80.34
$ws.Range("B36").Value = 6.53
$ws.Range("C36").Value = 37.76
$ws.Range("D36").Value = 23.73
$ws.Range("E36").Value = 12.06
$ws.Range("F36").Value = 10.06
$ws.Range("G36").Value = 9.029999999999999
$ws.Range("H36").Value = 59.05
$ws.Range("B37").Value = 5.05
$ws.Range("C37").Value = 20.68
$ws.Range("D37").Value = 9.5
$ws.Range("E37").Value = 7.35
$ws.Range("F37").Value = 6.94
$ws.Range("G37").Value = 8.74
$ws.Range("H37").Value = 39.65
$ws.Range("B38").Value = 5.51
$ws.Range("C38").Value = 23.18
$ws.Range("D38").Value = 14.57
$ws.Range("E38").Value = 9.35
$ws.Range("F38").Value = 9.02
$ws.Range("G38").Value = 8.029999999999999
$ws.Range("H38").Value = 32.38

# Remove the now-empty trailing rows (old rows 39-41: goiás/distrito federal duplicate labels + source footnote)
$ws.Range("A39:H41").ClearContents()

Write-Host "edit applied"
